$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tracker")
$ws.Range("A9").Value = 44971.939583333333
$ws.Range("B9").Value = "Submission"
$ws.Range("D9").Value = "LogReg"
$ws.Range("E9").Value = "scaling"
$ws.Range("H9").Value = 0.70799999999999996
$ws.Range("J9").Value = 0.42299999999999999
$ws.Range("K9").Value = "Maria"
$ws.Range("J10").Select()
